$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value cell replacements (1-indexed rows)
$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "1984"
$t.Cell(6, 1).Range.Text  = "0.00090"
$t.Cell(7, 1).Range.Text  = "0.00016"
$t.Cell(9, 1).Range.Text  = "0.00027"
$t.Cell(10, 1).Range.Text = "0.00030"
$t.Cell(11, 1).Range.Text = "0.00036"
$t.Cell(12, 1).Range.Text = "0.36347"

# Collapsed multi-run cells (tab-separated stat rows) -> single summary value
$t.Cell(44, 1).Range.Text = "99.62"
$t.Cell(45, 1).Range.Text = "0.36"
$t.Cell(46, 1).Range.Text = "95"
